# Fix bug in auto_run.sh
#
# The only real content change in this workbook edit is on the first
# worksheet ("C_N_New.conf"): cell B4 gains a "@FullAvg" suffix, turning
#   catalog=010A10;<SetCoord>(C_N, , subnodeid, id)->GEOMETRY
# into
#   catalog=010A10;<SetCoord>(C_N, , subnodeid, id)->GEOMETRY@FullAvg
#
# (Everything else in the raw XML diff -- fileVersion/rupBuild bumps, new
# xr/xr6/xr10/xr2 namespaces, absPath, revisionPtr, window position,
# selection cell, and the shuffled sharedStrings/cellXfs indices -- are
# artifacts Excel re-writes on every save and do not change any cell's
# visible value or effective formatting.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_N_New.conf")

$cell = $ws.Range("B4")
$cell.Value = $cell.Value + "@FullAvg"
